$wb = $excel.ActiveWorkbook

# Add the new sheet after the last existing sheet (HomePage), so it lands
# at the end of the tab strip and becomes the active sheet, matching the
# target workbook layout: Table, HomePage, RegistrationPage.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "RegistrationPage"

$ws.Range("A1").Value = "Gender"
$ws.Range("B1").Value = "Female"
$ws.Range("A2").Value = "First name"
$ws.Range("B2").Value = "First"
$ws.Range("A3").Value = "Last name"
$ws.Range("B3").Value = "Last"
$ws.Range("A4").Value = "Password"
# Write "Confirm password" before the "test123" values so the shared-string
# table fills in the same order as the target workbook (Password,
# Confirm password, test123).
$ws.Range("A5").Value = "Confirm password"
$ws.Range("B4").Value = "test123"
$ws.Range("B5").Value = "test123"

# Column A width to match the source registration sheet.
$ws.Columns.Item(1).ColumnWidth = 15.5

$ws.PageSetup.Orientation = 1

[void]$ws.Range("G11").Select()
